# fix cell string && add example && readme
#
# - C1 was "{C1}" -> now the "unmatched type" example "{unmatchedType}"
# - D1 was "{D1}" -> now the "undefined" example "{undefined}"
# - A2 was "hi {A2}" -> now "i am {A2}!"
# - C2 ("hello {C2} -") removed
# - New example added in B6: "{x} love {y}"
# - Column widths set (readme/example formatting)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- fix / replace example cell strings --
$ws.Range("C1").Value = "{unmatchedType}"
$ws.Range("D1").Value = "{undefined}"
$ws.Range("A2").Value = "i am {A2}!"
$ws.Range("C2").ClearContents()

# -- add new example row --
$ws.Range("B6").Value = "{x} love {y}"

# -- column widths (readme-style formatting) --
$ws.Columns.Item(1).ColumnWidth = 13.08
$ws.Columns.Item(2).ColumnWidth = 14.22
$ws.Columns.Item(3).ColumnWidth = 16.93
$ws.Columns.Item(4).ColumnWidth = 10.93

# -- move selection to the newly added example cell --
[void]$ws.Range("B6").Select()
